# "Alert Parameters (working).xlsx" - Third Iteration sheet:
# the "comments" column (Q) is repurposed into a sequential "message_ID"
# column, numbering every alert row 1-16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Third Iteration")

# Rename the column header from "comments" to "message_ID".
$ws.Range("Q1").Value = "message_ID"

# Number each data row (2-17) sequentially in column Q. Row 15 previously
# held a leftover text comment ("subtraction and alert value are in
# hours") - that note is replaced by its numeric message ID like every
# other row.
$lastRow = 17
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 17).Value = $row - 1
}

# Reflect the saved view state: scrolled right a bit with O23 selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$ws.Range("O23").Select()
